$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the header columns in row 2 (A2:I2) to match the new template order.
$newOrder = @(
    "Location_ID",
    "BusinessKey",
    "LocationTypeBusinessKey",
    "AreaKM",
    "Code",
    "Density",
    "Name",
    "ParentLocation_ID",
    "Population"
)

for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $newOrder[$i]
}
